$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "S k Constructions & Engineers"  ->  "S K Consultants & Engineers"
#    split across five runs: "S ", "K", " ", "Consultants", " & Engineers"
# ------------------------------------------------------------------
$findRange = $d.Content
$ok = $findRange.Find.Execute("S k Constructions & Engineers", $true, $false, `
                               $false, $false, $false, $true, 1, $false, "", 0)
if (-not $ok) {
    throw "Could not find target paragraph 'S k Constructions & Engineers'"
}
$start = $findRange.Start

# Original text layout (0-based offsets into the found range):
#   "S "(0-2)  "k"(2-3)  " "(3-4)  "Constructions"(4-17)  " & Engineers"(17-29)

# Fix the letter case and the word itself first (keeps everything in one run
# for the moment; we split it into separate runs afterwards).
$rK = $d.Range($start + 2, $start + 3)
$rK.Text = "K"

$rWord = $d.Range($start + 4, $start + 17)
$rWord.Text = "Consultants"

# New text layout after the (length-preserving / length-changing) edits above:
#   "S "(0-2)  "K"(2-3)  " "(3-4)  "Consultants"(4-15)  " & Engineers"(15-27)

# Force Word to keep these as distinct runs (even though their formatting is
# identical) by toggling a character property off and back on over each
# segment - this breaks run-merging at the segment boundaries.
$segments = @(
    @(0, 2),
    @(2, 3),
    @(3, 4),
    @(4, 15),
    @(15, 27)
)
foreach ($seg in $segments) {
    $r = $d.Range($start + $seg[0], $start + $seg[1])
    $r.Font.Bold = 0
    $r.Font.Bold = 1
}

# ------------------------------------------------------------------
# 2) Remove the whole paragraph that reads "18-BENG-ISM-822-43281707"
#    (including its paragraph mark, so the following paragraph takes its
#    place without leaving a blank line behind).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$targetIdx = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "18-BENG-ISM-822-43281707*") {
        $targetIdx = $i
    }
}
if ($targetIdx -eq -1) {
    throw "Could not find paragraph '18-BENG-ISM-822-43281707'"
}
$d.Paragraphs.Item($targetIdx).Range.Delete()
